# 自动更新Excel文件 - daily job that counts down the "剩余" (days remaining)
# column (E) by one for every shop row (2-99). Row 36 ("御盛园私房菜")
# carries a malformed start date (F36 = 202510929, not a real yyyymmdd
# value), so the real-world automation skips it instead of touching it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 99; $row++) {
    $dateCell = $ws.Cells.Item($row, 6)    # column F - start date (yyyymmdd)
    $dateText = [string]$dateCell.Text

    # A valid start date is an 8-digit yyyymmdd number; anything else
    # (blank, malformed, extra digits, ...) means this row is skipped.
    if ($dateText.Length -ne 8) {
        continue
    }

    $cell = $ws.Cells.Item($row, 5)        # column E - 剩余 (remaining)
    $current = $cell.Value2
    if ($null -ne $current) {
        $cell.Value = $current - 1
    }
}
